$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H10").Value = 39.5
$ws_ALC.Range("I10").Value = 39.5
$ws_ALC.Range("K10").Value = 39.5
$ws_ALC.Range("M10").Value = 253.5
$ws_ALC.Range("H12").Value = 436.66666
$ws_ALC.Range("I12").Value = 184
$ws_ALC.Range("J12").Value = 689.3333
$ws_ALC.Range("K12").Value = 184
$ws_ALC.Range("L12").Value = 689.3333
$ws_ALC.Range("M12").Value = -14
$ws_ALC.Range("N12").Value = -1029.3333
$ws_ALC.Range("H137").Value = 2758.4348
$ws_ALC.Range("I137").Value = 1268.5
$ws_ALC.Range("K137").Value = 3805.5
$ws_ALC.Range("M137").Value = -1255.5
$ws_ALC.Range("H138").Value = 4477.4165
$ws_ALC.Range("I138").Value = 4118.2
$ws_ALC.Range("J138").Value = 4734
$ws_ALC.Range("K138").Value = 12354.6
$ws_ALC.Range("L138").Value = 14202
$ws_ALC.Range("M138").Value = -7214.599999999999
$ws_ALC.Range("N138").Value = -24482
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("J10").Value = 0
$ws_ARM.Range("L10").Value = 0
$ws_ARM.Range("N10").Value = $null
$ws_ARM.Range("H11").Value = 3090
$ws_ARM.Range("I11").Value = 975
$ws_ARM.Range("J11").Value = 4500
$ws_ARM.Range("K11").Value = 975
$ws_ARM.Range("L11").Value = 4500
$ws_ARM.Range("M11").Value = -831
$ws_ARM.Range("N11").Value = -4788
$ws_ARM.Range("H92").Value = 52499.5
$ws_ARM.Range("J92").Value = 52499.5
$ws_ARM.Range("L92").Value = 52499.5
$ws_ARM.Range("N92").Value = -57491.5
$ws_ARM.Range("H132").Value = 2445.1177
$ws_ARM.Range("I132").Value = 2264.8333
$ws_ARM.Range("J132").Value = 2877.8
$ws_ARM.Range("K132").Value = 6794.499899999999
$ws_ARM.Range("L132").Value = 8633.400000000001
$ws_ARM.Range("M132").Value = -4264.499899999999
$ws_ARM.Range("N132").Value = -13693.4
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H76").Value = 15000
$ws_BSM.Range("J76").Value = 15000
$ws_BSM.Range("L76").Value = 15000
$ws_BSM.Range("N76").Value = -15630
$ws_BSM.Range("H79").Value = 15000
$ws_BSM.Range("J79").Value = 15000
$ws_BSM.Range("L79").Value = 15000
$ws_BSM.Range("N79").Value = -17184
$ws_BSM.Range("H105").Value = 2078.8
$ws_BSM.Range("I105").Value = 1848.5
$ws_BSM.Range("K105").Value = 1848.5
$ws_BSM.Range("M105").Value = -101.5
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H31").Value = 4308.4424
$ws_CRP.Range("I31").Value = 2395.4
$ws_CRP.Range("J31").Value = 4932.2607
$ws_CRP.Range("K31").Value = 2395.4
$ws_CRP.Range("L31").Value = 4932.2607
$ws_CRP.Range("M31").Value = -2100.4
$ws_CRP.Range("N31").Value = -5522.2607
$ws_CRP.Range("H34").Value = 4308.4424
$ws_CRP.Range("I34").Value = 2395.4
$ws_CRP.Range("J34").Value = 4932.2607
$ws_CRP.Range("K34").Value = 2395.4
$ws_CRP.Range("L34").Value = 4932.2607
$ws_CRP.Range("M34").Value = -2193.4
$ws_CRP.Range("N34").Value = -5336.2607
$ws_CRP.Range("H74").Value = 0
$ws_CRP.Range("I74").Value = 0
$ws_CRP.Range("J74").Value = 0
$ws_CRP.Range("K74").Value = 0
$ws_CRP.Range("L74").Value = $null
$ws_CRP.Range("M74").Value = $null
$ws_CRP.Range("N74").Value = 0
$ws_CRP.Range("H77").Value = 0
$ws_CRP.Range("I77").Value = 0
$ws_CRP.Range("J77").Value = 0
$ws_CRP.Range("K77").Value = 0
$ws_CRP.Range("L77").Value = $null
$ws_CRP.Range("M77").Value = $null
$ws_CRP.Range("N77").Value = 0
$ws_CRP.Range("H88").Value = 12500.25
$ws_CRP.Range("J88").Value = 12500.25
$ws_CRP.Range("L88").Value = 12500.25
$ws_CRP.Range("N88").Value = -13312.25
$ws_CRP.Range("H91").Value = 12500.25
$ws_CRP.Range("J91").Value = 12500.25
$ws_CRP.Range("L91").Value = 12500.25
$ws_CRP.Range("N91").Value = -15308.25
$ws_CRP.Range("H132").Value = 1623.25
$ws_CRP.Range("I132").Value = 1063.125
$ws_CRP.Range("K132").Value = 3189.375
$ws_CRP.Range("M132").Value = -659.375
$ws_CRP.Range("H134").Value = 1081.9565
$ws_CRP.Range("I134").Value = 1085.6818
$ws_CRP.Range("K134").Value = 3257.0454
$ws_CRP.Range("M134").Value = -722.0454
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H7").Value = 209.27586
$ws_CUL.Range("I7").Value = 98.210526
$ws_CUL.Range("K7").Value = 294.631578
$ws_CUL.Range("M7").Value = -182.631578
$ws_CUL.Range("H11").Value = 143027.72
$ws_CUL.Range("I11").Value = 333365
$ws_CUL.Range("J11").Value = 274.75
$ws_CUL.Range("K11").Value = 1000095
$ws_CUL.Range("L11").Value = 824.25
$ws_CUL.Range("M11").Value = -999955
$ws_CUL.Range("N11").Value = -1104.25
$ws_CUL.Range("H26").Value = 92
$ws_CUL.Range("I26").Value = 98.57143000000001
$ws_CUL.Range("K26").Value = 295.71429
$ws_CUL.Range("M26").Value = -7.714290000000005
$ws_CUL.Range("H32").Value = 160
$ws_CUL.Range("I32").Value = 160
$ws_CUL.Range("K32").Value = 480
$ws_CUL.Range("M32").Value = -197
$ws_CUL.Range("H34").Value = 1121.5294
$ws_CUL.Range("J34").Value = 1524.5
$ws_CUL.Range("L34").Value = 4573.5
$ws_CUL.Range("N34").Value = -4741.5
$ws_CUL.Range("H39").Value = 4812.7
$ws_CUL.Range("J39").Value = 5312.4443
$ws_CUL.Range("L39").Value = 15937.3329
$ws_CUL.Range("N39").Value = -16525.3329
$ws_CUL.Range("H48").Value = 7968.143
$ws_CUL.Range("J48").Value = 7968.143
$ws_CUL.Range("L48").Value = 23904.429
$ws_CUL.Range("N48").Value = -24404.429
$ws_CUL.Range("H55").Value = 5833.6924
$ws_CUL.Range("J55").Value = 6236.5
$ws_CUL.Range("L55").Value = 18709.5
$ws_CUL.Range("N55").Value = -19063.5
$ws_CUL.Range("H129").Value = 1798.2858
$ws_CUL.Range("I129").Value = 866.1667
$ws_CUL.Range("K129").Value = 2598.5001
$ws_CUL.Range("M129").Value = 2401.4999
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H4").Value = 0
$ws_GSM.Range("I4").Value = 0
$ws_GSM.Range("K4").Value = 0
$ws_GSM.Range("M4").Value = $null
$ws_GSM.Range("H13").Value = 68.333336
$ws_GSM.Range("I13").Value = 68.333336
$ws_GSM.Range("K13").Value = 68.333336
$ws_GSM.Range("M13").Value = 70.666664
$ws_GSM.Range("H40").Value = 24796.666
$ws_GSM.Range("I40").Value = 24790
$ws_GSM.Range("J40").Value = 24800
$ws_GSM.Range("K40").Value = 24790
$ws_GSM.Range("L40").Value = 24800
$ws_GSM.Range("M40").Value = -24639
$ws_GSM.Range("N40").Value = -25102
$ws_GSM.Range("H43").Value = 24310.111
$ws_GSM.Range("I43").Value = 99
$ws_GSM.Range("J43").Value = 31227.572
$ws_GSM.Range("K43").Value = 99
$ws_GSM.Range("L43").Value = 31227.572
$ws_GSM.Range("M43").Value = 52
$ws_GSM.Range("N43").Value = -31529.572
$ws_GSM.Range("H47").Value = 0
$ws_GSM.Range("J47").Value = 0
$ws_GSM.Range("L47").Value = $null
$ws_GSM.Range("N47").Value = 0
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H14").Value = 14726.394
$ws_LTW.Range("I14").Value = 14999
$ws_LTW.Range("J14").Value = 12750
$ws_LTW.Range("K14").Value = 14999
$ws_LTW.Range("L14").Value = 12750
$ws_LTW.Range("M14").Value = -14827
$ws_LTW.Range("N14").Value = -13094
$ws_LTW.Range("H16").Value = 750.8333
$ws_LTW.Range("I16").Value = 670
$ws_LTW.Range("K16").Value = 670
$ws_LTW.Range("M16").Value = -500
$ws_LTW.Range("H46").Value = 6754.3335
$ws_LTW.Range("I46").Value = 5263
$ws_LTW.Range("K46").Value = 5263
$ws_LTW.Range("M46").Value = -5075
$ws_LTW.Range("H98").Value = 64599.668
$ws_LTW.Range("J98").Value = 64599.668
$ws_LTW.Range("L98").Value = 64599.668
$ws_LTW.Range("N98").Value = -70589.66800000001
$ws_LTW.Range("H132").Value = 15165.833
$ws_LTW.Range("I132").Value = 9752
$ws_LTW.Range("K132").Value = 29256
$ws_LTW.Range("M132").Value = -26726
$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H103").Value = 31643.545
$ws_WVR.Range("J103").Value = 31643.545
$ws_WVR.Range("L103").Value = 31643.545
$ws_WVR.Range("N103").Value = -33987.545
$ws_WVR.Range("H132").Value = 1817.2667
$ws_WVR.Range("I132").Value = 1661.3572
$ws_WVR.Range("K132").Value = 4984.071599999999
$ws_WVR.Range("M132").Value = -2454.071599999999
$ws_WVR.Range("H136").Value = 3756.56
$ws_WVR.Range("I136").Value = 2988.7273
$ws_WVR.Range("J136").Value = 4359.857
$ws_WVR.Range("K136").Value = 8966.1819
$ws_WVR.Range("L136").Value = 13079.571
$ws_WVR.Range("M136").Value = -6416.1819
$ws_WVR.Range("N136").Value = -18179.571
